# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for rows 2-32 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 8
    3  = 7
    4  = 10
    5  = 8
    6  = 4
    7  = 8
    8  = 4
    9  = 5
    10 = 5
    11 = 1
    12 = 2
    13 = 0
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 7
    19 = 3
    20 = 5
    21 = 5
    22 = 4
    23 = 6
    24 = 8
    25 = 6
    26 = 6
    27 = 5
    28 = 5
    29 = 4
    30 = 5
    31 = 0
    32 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
